$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain TEXT even when it looks numeric,
# by temporarily forcing a Text number format, then clearing the format
# afterwards so the cell keeps its original (default) style/appearance.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "22.410.46"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.566.73"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.06%  "
Set-TextValue $ws.Range("D6") "284.78"
$ws.Range("E6").Value = "  -2.42%  "
Set-TextValue $ws.Range("D7") "0.3650"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  -2.71%  "
Set-TextValue $ws.Range("D9") "0.3329"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -1.84%  "
Set-TextValue $ws.Range("D11") "0.07396"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  -2.85%  "
Set-TextValue $ws.Range("D14") "5.949"
$ws.Range("E14").Value = "  -1.23%  "
Set-TextValue $ws.Range("D15") "6.899"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "1.566.14"
$ws.Range("E16").Value = "  -0.42%  "
Set-TextValue $ws.Range("D17") "0.00001103"
$ws.Range("E17").Value = "  -1.82%  "
Set-TextValue $ws.Range("D18") "88.06"
$ws.Range("E18").Value = "  -3.17%  "
Set-TextValue $ws.Range("D19") "0.06699"
$ws.Range("E19").Value = "  -0.78%  "
Set-TextValue $ws.Range("D21") "6.350"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("E22").Value = "  -0.93%  "
Set-TextValue $ws.Range("D23") "12.01"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "22.411.76"
$ws.Range("E24").Value = "  -0.22%  "
Set-TextValue $ws.Range("D25") "2.386"
$ws.Range("E25").Value = "  +1.84%  "
Set-TextValue $ws.Range("D26") "2.546"
$ws.Range("E26").Value = "  -5.37%  "
Set-TextValue $ws.Range("D27") "150.14"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  -3.72%  "
Set-TextValue $ws.Range("D29") "4.995"
$ws.Range("E29").Value = "  -0.20%  "
Set-TextValue $ws.Range("D30") "123.77"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "1.742.81"
$ws.Range("E31").Value = "  -0.38%  "
Set-TextValue $ws.Range("D32") "1.053"
$ws.Range("E32").Value = "  -0.66%  "
Set-TextValue $ws.Range("D33") "6.102"
$ws.Range("E33").Value = "  -1.70%  "
Set-TextValue $ws.Range("D34") "1.989"
$ws.Range("E34").Value = "  +0.12%  "
Set-TextValue $ws.Range("D35") "9.794"
$ws.Range("E35").Value = "  -0.42%  "
Set-TextValue $ws.Range("D36") "0.08274"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  -2.92%  "
Set-TextValue $ws.Range("D39") "0.06412"
$ws.Range("E39").Value = "  -1.75%  "
Set-TextValue $ws.Range("D40") "5.365"
$ws.Range("E40").Value = "  -2.23%  "
Set-TextValue $ws.Range("D41") "1.280"
$ws.Range("E41").Value = "  -7.00%  "
Set-TextValue $ws.Range("D42") "0.6247"
$ws.Range("E42").Value = "  +0.18%  "
Set-TextValue $ws.Range("D43") "11.15"
$ws.Range("E43").Value = "  -1.41%  "
Set-TextValue $ws.Range("D45") "13.79"
$ws.Range("E45").Value = "  -1.24%  "
Set-TextValue $ws.Range("D46") "0.6038"
$ws.Range("E46").Value = "  +3.80%  "
Set-TextValue $ws.Range("D47") "3.748"
$ws.Range("E47").Value = "  -1.75%  "
Set-TextValue $ws.Range("D48") "2.032"
$ws.Range("E48").Value = "  -2.27%  "
Set-TextValue $ws.Range("D49") "123.88"
$ws.Range("E49").Value = "  -4.49%  "
Set-TextValue $ws.Range("D50") "1.213"
Set-TextValue $ws.Range("D51") "0.07200"
$ws.Range("E51").Value = "  -1.80%  "
